$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rooting-decision columns (C: Root Outgroups?, D: Root Ictiobinae,
#     E: Root Moxostomatini, F: Root Catostominae, G: Notes) for loci
#     L119-L175 (rows 113-167) ---
$ws.Range("C113").Value = 'N'
$ws.Range("D113").Value = 'N'
$ws.Range("E113").Value = 'N'
$ws.Range("F113").Value = 'Y'
$ws.Range("G113").Value = 'Makes deep bodied clade'
$ws.Range("C114").Value = 'N'
$ws.Range("D114").Value = 'Y'
$ws.Range("G114").Value = 'Catostomine tribes make sense'
$ws.Range("C115").Value = 'N'
$ws.Range("D115").Value = 'Y'
$ws.Range("G115").Value = 'Only Ictiobus niger present'
$ws.Range("C116").Value = 'N'
$ws.Range("D116").Value = 'N'
$ws.Range("E116").Value = 'N'
$ws.Range("F116").Value = 'N'
$ws.Range("G116").Value = 'Useless'
$ws.Range("C117").Value = 'Y'
$ws.Range("C119").Value = 'N'
$ws.Range("D119").Value = 'Y'
$ws.Range("G119").Value = 'Moxostomatini and Catostomini mixed.'
$ws.Range("C120").Value = 'N'
$ws.Range("D120").Value = 'N'
$ws.Range("E120").Value = 'N'
$ws.Range("F120").Value = 'N'
$ws.Range("G120").Value = 'Useless'
$ws.Range("C124").Value = 'Y'
$ws.Range("C125").Value = 'Y'
$ws.Range("C127").Value = 'Y'
$ws.Range("C134").Value = 'Y'
$ws.Range("C135").Value = 'N'
$ws.Range("D135").Value = 'N'
$ws.Range("E135").Value = 'N'
$ws.Range("F135").Value = 'Y'
$ws.Range("G135").Value = 'Makes deep bodied clade'
$ws.Range("C136").Value = 'N'
$ws.Range("D136").Value = 'Y'
$ws.Range("G136").Value = 'Moxostomatini and Catostomini mixed.'
$ws.Range("C137").Value = 'Y'
$ws.Range("C139").Value = 'N'
$ws.Range("D139").Value = 'Y'
$ws.Range("G139").Value = 'Catostominae makes sense.'
$ws.Range("C140").Value = 'N'
$ws.Range("D140").Value = 'Y'
$ws.Range("G140").Value = 'Moxostoma is polyphyletic.'
$ws.Range("C141").Value = 'N'
$ws.Range("D141").Value = 'Y'
$ws.Range("G141").Value = 'Makes a deep-bodied clade.'
$ws.Range("C142").Value = 'Y'
$ws.Range("C143").Value = 'N'
$ws.Range("D143").Value = 'Y'
$ws.Range("G143").Value = 'Catostominae is weird; Deep bodied clade formed.'
$ws.Range("C144").Value = 'N'
$ws.Range("D144").Value = 'Y'
$ws.Range("G144").Value = 'Catostominae is weird'
$ws.Range("C145").Value = 'Y'
$ws.Range("C146").Value = 'Y'
$ws.Range("C149").Value = 'N'
$ws.Range("D149").Value = 'N'
$ws.Range("E149").Value = 'N'
$ws.Range("F149").Value = 'Y'
$ws.Range("G149").Value = 'Catostominae is weird.'
$ws.Range("C151").Value = 'Y'
$ws.Range("C152").Value = 'Y'
$ws.Range("C154").Value = 'N'
$ws.Range("D154").Value = 'Y'
$ws.Range("C155").Value = 'Y'
$ws.Range("C157").Value = 'Y'
$ws.Range("C159").Value = 'N'
$ws.Range("D159").Value = 'N'
$ws.Range("E159").Value = 'N'
$ws.Range("F159").Value = 'Y'
$ws.Range("G159").Value = 'Makes deep bodied clade'
$ws.Range("C163").Value = 'N'
$ws.Range("D163").Value = 'N'
$ws.Range("E163").Value = 'N'
$ws.Range("F163").Value = 'Y'
$ws.Range("G163").Value = 'Makes deep bodied clade'
$ws.Range("C164").Value = 'Y'
$ws.Range("C165").Value = 'Y'
$ws.Range("C166").Value = 'Y'
$ws.Range("C167").Value = 'Y'

# --- Summary rows 269 (per-column Y counts) and 270 (grand total) ---
$ws.Range("C269").Formula = '=COUNTIF(C2:C268, "Y")'
$ws.Range("D269:F269").Formula = '=COUNTIF(D2:D268, "Y")'

$ws.Range("A270").Value = "SUM->"
$ws.Range("A270").Font.Bold = $true
$ws.Range("B270").Formula = "=SUM(C269:F269)"

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Restore the view: scroll/select back near the data that was edited ---
$ws.Range("A158").Select()
$ws.Range("C167").Select()
